# Update "detection field test data.xlsx" - add two new field-test rows
# (28 May 2025 / ABSENCE search, and 30 May 2025 / PRESENCE search) to the
# "dog" sheet, and move the selection down to where the user was last
# working (L16, scrolled so column E is at the left edge).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dog")
$ws.Activate()

# --- Row 14: 28 May 2025, ABSENCE search ------------------------------
# Borrow the number formats from row 7, which already has the same
# text/date/time column layout (Type + "NA" placed-time text + search
# time) that this new row needs.
$ws.Range("A7:L7").Copy() | Out-Null
$ws.Range("A14:L14").PasteSpecial(-4122) | Out-Null

$ws.Range("A14").Value = 45805
$ws.Range("B14").Value = "ABSENCE"
$ws.Range("C14").Value = "NA"
$ws.Range("D14").Value = 0.5625
$ws.Range("E14").Value = 15
$ws.Range("F14").Value = 15
$ws.Range("G14").Value = "Sunny, cool"
$ws.Range("H14").Value = $false
$ws.Range("I14").Value = "14 minutes 9 seconds"
$ws.Range("J14").Value = 849
$ws.Range("K14").Value = "NA"
$ws.Range("L14").Value = "Worked downhill. Koda did a really thorough sweep so I was fairly convinced there was no target, but we still completed the zig zag transects before finishing the search."

# --- Row 15: 30 May 2025, PRESENCE search -------------------------------
# Borrow the number formats from row 10, which has numeric "time target
# placed"/"time of search" columns like this new row.
$ws.Range("A10:L10").Copy() | Out-Null
$ws.Range("A15:L15").PasteSpecial(-4122) | Out-Null

$ws.Range("A15").Value = 45807
$ws.Range("B15").Value = "PRESENCE"
$ws.Range("C15").Value = 0.29166666666666669
$ws.Range("D15").Value = 0.52777777777777779
$ws.Range("E15").Value = 16
$ws.Range("F15").Value = 13
$ws.Range("G15").Value = "Sunny, mild"
$ws.Range("H15").Value = $true
$ws.Range("I15").Value = "1 minute 12 seconds"
$ws.Range("J15").Value = 72
$ws.Range("K15").Value = "Primary sweeps"
$ws.Range("L15").Value = "Sent Koda in and she searched downhill and got onto odour in less than a minute."

$excel.CutCopyMode = $false

# --- View state: scroll / select where the user left off ---------------
$ws.Range("E1").Select() | Out-Null
$excel.ActiveWindow.ScrollColumn = 5
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("L16").Select() | Out-Null
